$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated figures) ---
$ws.Range("G10").Value = 0.185551295138736
$ws.Range("G11").Value = 0.185551295138736
$ws.Range("G23").Value = 0.191440605331555
$ws.Range("G24").Value = 0.191440605331555
$ws.Range("G36").Value = 1325.20241416136
$ws.Range("H36").Value = 12939.3521333977
$ws.Range("G37").Value = 1325.20241416136
$ws.Range("H37").Value = 12939.3521333977
$ws.Range("G38").Value = 1325.20241416136
$ws.Range("H38").Value = 12939.3521333977
$ws.Range("G39").Value = 1325.20241416136
$ws.Range("H39").Value = 12939.3521333977
$ws.Range("G40").Value = 0.212893082006662
$ws.Range("G41").Value = 0.212893082006662
$ws.Range("G53").Value = 1192.21075214396
$ws.Range("H53").Value = 12939.3521333977
$ws.Range("G54").Value = 1192.21075214396
$ws.Range("H54").Value = 12939.3521333977
$ws.Range("G55").Value = 1192.21075214396
$ws.Range("H55").Value = 12939.3521333977
$ws.Range("G56").Value = 1192.21075214396
$ws.Range("H56").Value = 12939.3521333977
$ws.Range("G57").Value = 0.205453303968608
$ws.Range("G58").Value = 0.205453303968608
$ws.Range("G70").Value = 1152.80286681662
$ws.Range("H70").Value = 12939.3521333977
$ws.Range("G71").Value = 1152.80286681662
$ws.Range("H71").Value = 12939.3521333977
$ws.Range("G72").Value = 1152.80286681662
$ws.Range("H72").Value = 12939.3521333977
$ws.Range("G73").Value = 1152.80286681662
$ws.Range("H73").Value = 12939.3521333977
$ws.Range("G74").Value = 0.214416548499924
$ws.Range("G75").Value = 0.214416548499924
$ws.Range("G87").Value = 1103.76120014995
$ws.Range("H87").Value = 12939.3521333977
$ws.Range("G88").Value = 1103.76120014995
$ws.Range("H88").Value = 12939.3521333977
$ws.Range("G89").Value = 1103.76120014995
$ws.Range("H89").Value = 12939.3521333977
$ws.Range("G90").Value = 1103.76120014995
$ws.Range("H90").Value = 12939.3521333977
$ws.Range("G104").Value = 1588.01795230982
$ws.Range("H104").Value = 12939.3521333977
$ws.Range("I104").Value = 8152.13241
$ws.Range("G105").Value = 1588.01795230982
$ws.Range("H105").Value = 12939.3521333977
$ws.Range("I105").Value = 8152.13241
$ws.Range("G106").Value = 1588.01795230982
$ws.Range("H106").Value = 12939.3521333977
$ws.Range("I106").Value = 8152.13241
$ws.Range("G107").Value = 1588.01795230982
$ws.Range("H107").Value = 12939.3521333977
$ws.Range("I107").Value = 8152.13241

# --- Append new rows 118-134 ---
# Row 118
$ws.Range("A118").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B118").Value = "Visual Clarity (Sediment class 1)"
$ws.Range("C118").Value = "D"
$ws.Range("D118").Value = "2019 - 2023"
$ws.Range("E118").Value = "RepSite"
$ws.Range("F118").Value = 0.305
$ws.Range("G118").Value = 0.2975
$ws.Range("H118").Value = 0.48
$ws.Range("I118").Value = 0.478
$ws.Range("J118").Value = ""
$ws.Range("K118").Value = ""
$ws.Range("L118").Value = 0.25
$ws.Range("M118").Value = 0.3972
$ws.Range("N118").Value = 0.46626
$ws.Range("O118").Value = 1791973.7
$ws.Range("P118").Value = 5504665
$ws.Range("Q118").Value = "Horowhenua District"
$ws.Range("R118").Value = "Waiopehu"
$ws.Range("S118").Value = "Lake Horowhenua"
$ws.Range("T118").Value = "Hoki_1a"
$ws.Range("U118").Value = "m"

# Row 119
$ws.Range("A119").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B119").Value = "DRP (95th Percentile)"
$ws.Range("C119").Value = "D"
$ws.Range("D119").Value = "2019 - 2023"
$ws.Range("E119").Value = "RepSite"
$ws.Range("F119").Value = 0.03
$ws.Range("G119").Value = 0.03378
$ws.Range("H119").Value = 0.08699999999999999
$ws.Range("I119").Value = 0.062
$ws.Range("J119").Value = ""
$ws.Range("K119").Value = ""
$ws.Range("L119").Value = 0.036
$ws.Range("M119").Value = 0.048
$ws.Range("N119").Value = 0.0595
$ws.Range("O119").Value = 1791973.7
$ws.Range("P119").Value = 5504665
$ws.Range("Q119").Value = "Horowhenua District"
$ws.Range("R119").Value = "Waiopehu"
$ws.Range("S119").Value = "Lake Horowhenua"
$ws.Range("T119").Value = "Hoki_1a"
$ws.Range("U119").Value = "mg/L"

# Row 120
$ws.Range("A120").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B120").Value = "DRP (Median)"
$ws.Range("C120").Value = "D"
$ws.Range("D120").Value = "2019 - 2023"
$ws.Range("E120").Value = "RepSite"
$ws.Range("F120").Value = 0.03
$ws.Range("G120").Value = 0.03378
$ws.Range("H120").Value = 0.08699999999999999
$ws.Range("I120").Value = 0.062
$ws.Range("J120").Value = ""
$ws.Range("K120").Value = ""
$ws.Range("L120").Value = 0.036
$ws.Range("M120").Value = 0.048
$ws.Range("N120").Value = 0.0595
$ws.Range("O120").Value = 1791973.7
$ws.Range("P120").Value = 5504665
$ws.Range("Q120").Value = "Horowhenua District"
$ws.Range("R120").Value = "Waiopehu"
$ws.Range("S120").Value = "Lake Horowhenua"
$ws.Range("T120").Value = "Hoki_1a"
$ws.Range("U120").Value = "mg/L"

# Row 121
$ws.Range("A121").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B121").Value = "E coli (>260)"
$ws.Range("C121").Value = "E"
$ws.Range("D121").Value = "2019 - 2023"
$ws.Range("E121").Value = "RepSite"
$ws.Range("F121").Value = 355
$ws.Range("G121").Value = 1089.76648207347
$ws.Range("H121").Value = 10421.3241036736
$ws.Range("I121").Value = 5199
$ws.Range("J121").Value = 34
$ws.Range("K121").Value = 58
$ws.Range("L121").Value = 568
$ws.Range("M121").Value = 1600
$ws.Range("N121").Value = 4042.5
$ws.Range("O121").Value = 1791973.7
$ws.Range("P121").Value = 5504665
$ws.Range("Q121").Value = "Horowhenua District"
$ws.Range("R121").Value = "Waiopehu"
$ws.Range("S121").Value = "Lake Horowhenua"
$ws.Range("T121").Value = "Hoki_1a"
$ws.Range("U121").Value = "% exceedances over 260/100 mL"

# Row 122
$ws.Range("A122").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B122").Value = "E coli (>540)"
$ws.Range("C122").Value = "E"
$ws.Range("D122").Value = "2019 - 2023"
$ws.Range("E122").Value = "RepSite"
$ws.Range("F122").Value = 355
$ws.Range("G122").Value = 1089.76648207347
$ws.Range("H122").Value = 10421.3241036736
$ws.Range("I122").Value = 5199
$ws.Range("J122").Value = 34
$ws.Range("K122").Value = 58
$ws.Range("L122").Value = 568
$ws.Range("M122").Value = 1600
$ws.Range("N122").Value = 4042.5
$ws.Range("O122").Value = 1791973.7
$ws.Range("P122").Value = 5504665
$ws.Range("Q122").Value = "Horowhenua District"
$ws.Range("R122").Value = "Waiopehu"
$ws.Range("S122").Value = "Lake Horowhenua"
$ws.Range("T122").Value = "Hoki_1a"
$ws.Range("U122").Value = "% exceedances over 540/100 mL"

# Row 123
$ws.Range("A123").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B123").Value = "E coli (Median)"
$ws.Range("C123").Value = "E"
$ws.Range("D123").Value = "2019 - 2023"
$ws.Range("E123").Value = "RepSite"
$ws.Range("F123").Value = 355
$ws.Range("G123").Value = 1089.76648207347
$ws.Range("H123").Value = 10421.3241036736
$ws.Range("I123").Value = 5199
$ws.Range("J123").Value = 34
$ws.Range("K123").Value = 58
$ws.Range("L123").Value = 568
$ws.Range("M123").Value = 1600
$ws.Range("N123").Value = 4042.5
$ws.Range("O123").Value = 1791973.7
$ws.Range("P123").Value = 5504665
$ws.Range("Q123").Value = "Horowhenua District"
$ws.Range("R123").Value = "Waiopehu"
$ws.Range("S123").Value = "Lake Horowhenua"
$ws.Range("T123").Value = "Hoki_1a"
$ws.Range("U123").Value = "E. coli/100 mL"

# Row 124
$ws.Range("A124").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B124").Value = "E coli (95th Percentile)"
$ws.Range("C124").Value = "E"
$ws.Range("D124").Value = "2019 - 2023"
$ws.Range("E124").Value = "RepSite"
$ws.Range("F124").Value = 355
$ws.Range("G124").Value = 1089.76648207347
$ws.Range("H124").Value = 10421.3241036736
$ws.Range("I124").Value = 5199
$ws.Range("J124").Value = 34
$ws.Range("K124").Value = 58
$ws.Range("L124").Value = 568
$ws.Range("M124").Value = 1600
$ws.Range("N124").Value = 4042.5
$ws.Range("O124").Value = 1791973.7
$ws.Range("P124").Value = 5504665
$ws.Range("Q124").Value = "Horowhenua District"
$ws.Range("R124").Value = "Waiopehu"
$ws.Range("S124").Value = "Lake Horowhenua"
$ws.Range("T124").Value = "Hoki_1a"
$ws.Range("U124").Value = "E. coli/100 mL"

# Row 125
$ws.Range("A125").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B125").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C125").Value = "B"
$ws.Range("D125").Value = "2019 - 2023"
$ws.Range("E125").Value = "RepSite"
$ws.Range("F125").Value = 0.12488
$ws.Range("G125").Value = 0.135547826159867
$ws.Range("H125").Value = 0.41958041958042
$ws.Range("I125").Value = 0.32755
$ws.Range("J125").Value = ""
$ws.Range("K125").Value = ""
$ws.Range("L125").Value = 0.11194
$ws.Range("M125").Value = 0.25351
$ws.Range("N125").Value = 0.3135
$ws.Range("O125").Value = 1791973.7
$ws.Range("P125").Value = 5504665
$ws.Range("Q125").Value = "Horowhenua District"
$ws.Range("R125").Value = "Waiopehu"
$ws.Range("S125").Value = "Lake Horowhenua"
$ws.Range("T125").Value = "Hoki_1a"
$ws.Range("U125").Value = "mg NH4-N/L"

# Row 126
$ws.Range("A126").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B126").Value = "Ammoniacal-N (Median)"
$ws.Range("C126").Value = "B"
$ws.Range("D126").Value = "2019 - 2023"
$ws.Range("E126").Value = "RepSite"
$ws.Range("F126").Value = 0.12488
$ws.Range("G126").Value = 0.135547826159867
$ws.Range("H126").Value = 0.41958041958042
$ws.Range("I126").Value = 0.32755
$ws.Range("J126").Value = ""
$ws.Range("K126").Value = ""
$ws.Range("L126").Value = 0.11194
$ws.Range("M126").Value = 0.25351
$ws.Range("N126").Value = 0.3135
$ws.Range("O126").Value = 1791973.7
$ws.Range("P126").Value = 5504665
$ws.Range("Q126").Value = "Horowhenua District"
$ws.Range("R126").Value = "Waiopehu"
$ws.Range("S126").Value = "Lake Horowhenua"
$ws.Range("T126").Value = "Hoki_1a"
$ws.Range("U126").Value = "mg NH4-N/L"

# Row 127
$ws.Range("A127").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B127").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C127").Value = "B"
$ws.Range("D127").Value = "2019 - 2023"
$ws.Range("E127").Value = "RepSite"
$ws.Range("F127").Value = 1.825
$ws.Range("G127").Value = 1.76746
$ws.Range("H127").Value = 3.07
$ws.Range("I127").Value = 2.71
$ws.Range("J127").Value = ""
$ws.Range("K127").Value = ""
$ws.Range("L127").Value = 1.17
$ws.Range("M127").Value = 2.35
$ws.Range("N127").Value = 2.65
$ws.Range("O127").Value = 1791973.7
$ws.Range("P127").Value = 5504665
$ws.Range("Q127").Value = "Horowhenua District"
$ws.Range("R127").Value = "Waiopehu"
$ws.Range("S127").Value = "Lake Horowhenua"
$ws.Range("T127").Value = "Hoki_1a"
$ws.Range("U127").Value = "mg NO3-N/L"

# Row 128
$ws.Range("A128").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B128").Value = "Nitrate-N (Median)"
$ws.Range("C128").Value = "B"
$ws.Range("D128").Value = "2019 - 2023"
$ws.Range("E128").Value = "RepSite"
$ws.Range("F128").Value = 1.825
$ws.Range("G128").Value = 1.76746
$ws.Range("H128").Value = 3.07
$ws.Range("I128").Value = 2.71
$ws.Range("J128").Value = ""
$ws.Range("K128").Value = ""
$ws.Range("L128").Value = 1.17
$ws.Range("M128").Value = 2.35
$ws.Range("N128").Value = 2.65
$ws.Range("O128").Value = 1791973.7
$ws.Range("P128").Value = 5504665
$ws.Range("Q128").Value = "Horowhenua District"
$ws.Range("R128").Value = "Waiopehu"
$ws.Range("S128").Value = "Lake Horowhenua"
$ws.Range("T128").Value = "Hoki_1a"
$ws.Range("U128").Value = "mg NO3-N/L"

# Row 129
$ws.Range("A129").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B129").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("C129").Value = ""
$ws.Range("D129").Value = "2019 - 2023"
$ws.Range("E129").Value = "RepSite"
$ws.Range("F129").Value = 2.189
$ws.Range("G129").Value = 2.18084
$ws.Range("H129").Value = 4.064
$ws.Range("I129").Value = 3.176
$ws.Range("J129").Value = ""
$ws.Range("K129").Value = ""
$ws.Range("L129").Value = 2.006
$ws.Range("M129").Value = 2.844
$ws.Range("N129").Value = 3.042
$ws.Range("O129").Value = 1791973.7
$ws.Range("P129").Value = 5504665
$ws.Range("Q129").Value = "Horowhenua District"
$ws.Range("R129").Value = "Waiopehu"
$ws.Range("S129").Value = "Lake Horowhenua"
$ws.Range("T129").Value = "Hoki_1a"
$ws.Range("U129").Value = "g/m3"

# Row 130
$ws.Range("A130").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B130").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("C130").Value = ""
$ws.Range("D130").Value = "2019 - 2023"
$ws.Range("E130").Value = "RepSite"
$ws.Range("F130").Value = 2.189
$ws.Range("G130").Value = 2.18084
$ws.Range("H130").Value = 4.064
$ws.Range("I130").Value = 3.176
$ws.Range("J130").Value = ""
$ws.Range("K130").Value = ""
$ws.Range("L130").Value = 2.006
$ws.Range("M130").Value = 2.844
$ws.Range("N130").Value = 3.042
$ws.Range("O130").Value = 1791973.7
$ws.Range("P130").Value = 5504665
$ws.Range("Q130").Value = "Horowhenua District"
$ws.Range("R130").Value = "Waiopehu"
$ws.Range("S130").Value = "Lake Horowhenua"
$ws.Range("T130").Value = "Hoki_1a"
$ws.Range("U130").Value = "g/m3"

# Row 131
$ws.Range("A131").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B131").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("C131").Value = ""
$ws.Range("D131").Value = "2019 - 2023"
$ws.Range("E131").Value = "RepSite"
$ws.Range("F131").Value = 3.415
$ws.Range("G131").Value = 3.3394
$ws.Range("H131").Value = 6.25
$ws.Range("I131").Value = 4.4
$ws.Range("J131").Value = ""
$ws.Range("K131").Value = ""
$ws.Range("L131").Value = 2.99
$ws.Range("M131").Value = 3.97
$ws.Range("N131").Value = 4.325
$ws.Range("O131").Value = 1791973.7
$ws.Range("P131").Value = 5504665
$ws.Range("Q131").Value = "Horowhenua District"
$ws.Range("R131").Value = "Waiopehu"
$ws.Range("S131").Value = "Lake Horowhenua"
$ws.Range("T131").Value = "Hoki_1a"
$ws.Range("U131").Value = "g/m3"

# Row 132
$ws.Range("A132").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B132").Value = "Total Nitrogen (Median)"
$ws.Range("C132").Value = ""
$ws.Range("D132").Value = "2019 - 2023"
$ws.Range("E132").Value = "RepSite"
$ws.Range("F132").Value = 3.415
$ws.Range("G132").Value = 3.3394
$ws.Range("H132").Value = 6.25
$ws.Range("I132").Value = 4.4
$ws.Range("J132").Value = ""
$ws.Range("K132").Value = ""
$ws.Range("L132").Value = 2.99
$ws.Range("M132").Value = 3.97
$ws.Range("N132").Value = 4.325
$ws.Range("O132").Value = 1791973.7
$ws.Range("P132").Value = 5504665
$ws.Range("Q132").Value = "Horowhenua District"
$ws.Range("R132").Value = "Waiopehu"
$ws.Range("S132").Value = "Lake Horowhenua"
$ws.Range("T132").Value = "Hoki_1a"
$ws.Range("U132").Value = "g/m3"

# Row 133
$ws.Range("A133").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B133").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("C133").Value = ""
$ws.Range("D133").Value = "2019 - 2023"
$ws.Range("E133").Value = "RepSite"
$ws.Range("F133").Value = 0.1215
$ws.Range("G133").Value = 0.15678
$ws.Range("H133").Value = 0.776
$ws.Range("I133").Value = 0.341
$ws.Range("J133").Value = ""
$ws.Range("K133").Value = ""
$ws.Range("L133").Value = 0.128
$ws.Range("M133").Value = 0.224
$ws.Range("N133").Value = 0.309
$ws.Range("O133").Value = 1791973.7
$ws.Range("P133").Value = 5504665
$ws.Range("Q133").Value = "Horowhenua District"
$ws.Range("R133").Value = "Waiopehu"
$ws.Range("S133").Value = "Lake Horowhenua"
$ws.Range("T133").Value = "Hoki_1a"
$ws.Range("U133").Value = "g/m3"

# Row 134
$ws.Range("A134").Value = "L Horowhenua Inflow at Lindsay Road"
$ws.Range("B134").Value = "Total Phosphorus (Median)"
$ws.Range("C134").Value = ""
$ws.Range("D134").Value = "2019 - 2023"
$ws.Range("E134").Value = "RepSite"
$ws.Range("F134").Value = 0.1215
$ws.Range("G134").Value = 0.15678
$ws.Range("H134").Value = 0.776
$ws.Range("I134").Value = 0.341
$ws.Range("J134").Value = ""
$ws.Range("K134").Value = ""
$ws.Range("L134").Value = 0.128
$ws.Range("M134").Value = 0.224
$ws.Range("N134").Value = 0.309
$ws.Range("O134").Value = 1791973.7
$ws.Range("P134").Value = 5504665
$ws.Range("Q134").Value = "Horowhenua District"
$ws.Range("R134").Value = "Waiopehu"
$ws.Range("S134").Value = "Lake Horowhenua"
$ws.Range("T134").Value = "Hoki_1a"
$ws.Range("U134").Value = "g/m3"

